$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E..T for data rows 2..10 (row index aligns with Excel row number)
$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$newValues = @{}
$newValues[2] = @(3, 1, 0.29347, 0.8804099999999999, 0.1501202107524681, 0.1501202107524681, 2, 0.6666666666666666, 0.04110666666666667, 0.12332, 0.004039063627665411, 0.004039063627665412, 0.01206357346666666, 0.1085721612, 0.00060634508302776, 0.0006063450830277601)
$newValues[3] = @(3, 1, 0.29347, 0.8804099999999999, 0.1501202107524681, 0.1501202107524681, 3, 1, 7.100470666666666, 21.301412, 0.6976788714492014, 0.6976788714492014, 2.083775126546666, 18.75397613892, 0.1047356992194982, 0.1047356992194982)
$newValues[4] = @(3, 1, 0.29347, 0.8804099999999999, 0.1501202107524681, 0.1501202107524681, 3, 1, 3.035699, 9.107097, 0.2982820649231331, 0.2982820649231332, 0.8908865855299998, 8.017979269769999, 0.04477816644994213, 0.04477816644994214)
$newValues[5] = @(3, 1, 0.4262446666666667, 1.278734, 0.2180391153852712, 0.2180391153852712, 2, 0.6666666666666666, 0.04110666666666667, 0.12332, 0.004039063627665411, 0.004039063627665412, 0.01752149743111111, 0.15769347688, 0.0008806738603609906, 0.0008806738603609908)
$newValues[6] = @(3, 1, 0.4262446666666667, 1.278734, 0.2180391153852712, 0.2180391153852712, 3, 1, 7.100470666666666, 21.301412, 0.6976788714492014, 0.6976788714492014, 3.026537752489777, 27.238839772408, 0.1521212839537782, 0.1521212839537782)
$newValues[7] = @(3, 1, 0.4262446666666667, 1.278734, 0.2180391153852712, 0.2180391153852712, 3, 1, 3.035699, 9.107097, 0.2982820649231331, 0.2982820649231332, 1.293950508355333, 11.645554575198, 0.06503715757113199, 0.065037157571132)
$newValues[8] = @(3, 1, 1.235185333333333, 3.705556, 0.6318406738622607, 0.6318406738622606, 2, 0.6666666666666666, 0.04110666666666667, 0.12332, 0.004039063627665411, 0.004039063627665412, 0.05077435176888889, 0.45696916592, 0.00255204468427666, 0.002552044684276661)
$newValues[9] = @(3, 1, 1.235185333333333, 3.705556, 0.6318406738622607, 0.6318406738622606, 3, 1, 7.100470666666666, 21.301412, 0.6976788714492014, 0.6976788714492014, 8.770397227230221, 78.933575045072, 0.440821888275925, 0.4408218882759249)
$newValues[10] = @(3, 1, 1.235185333333333, 3.705556, 0.6318406738622607, 0.6318406738622606, 3, 1, 3.035699, 9.107097, 0.2982820649231331, 0.2982820649231332, 3.749650881214666, 33.746857930932, 0.188466740902059, 0.1884667409020591)

foreach ($r in $newValues.Keys) {
    $rowVals = $newValues[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$r"
        $ws.Range($addr).Value = $rowVals[$i]
    }
}

Write-Host "done"